# Applies the "search history" / TA-questions edit described in the commit.
$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Paragraph 1 ("Hörður" line): drop the "--lookið-display -- síðan function--"
# text but keep the preceding <w:tab/>.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("--lookið-display -- síðan function--", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Delete()
}

# Paragraph ("read only review  back - pay"): remove the whole paragraph
# (Expand to wdParagraph = 4 so the trailing paragraph mark is included too,
# which merges the now-empty paragraph away instead of leaving a blank one).
$rng = $d.Content
$found = $rng.Find.Execute("read only review  back - pay", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4) | Out-Null
    $rng.Delete()
}

# --- Change 2 -------------------------------------------------------------
# Remove the whole "Search/browse history -- sambærilegt og amazon" paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("Search/browse history -- sambærilegt og amazon", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4) | Out-Null
    $rng.Delete()
}

# Remove the whole "checkout (contact/payment info) ..." paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("checkout (contact/payment info) -- ekki hægt að ýta pay ef það vantar eitthvað  --- back - pay", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4) | Out-Null
    $rng.Delete()
}

# Move the "_GoBack" bookmark from its own (now otherwise empty) paragraph to
# the start of the "working checkout that keeps the information" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$rng = $d.Content
$found = $rng.Find.Execute("working checkout that keeps the information", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $startRng = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $startRng)
}
